$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8201449513435364
$ws.Range("B1").Value = 0.6396301984786987
$ws.Range("C1").Value = 4.707525253295898
$ws.Range("D1").Value = 2.798688411712646
$ws.Range("E1").Value = 1.199601054191589
